$d = $word.ActiveDocument

# --- Change 1: merge the split "disegnare la curva di Bézier..." runs into one run ---
$d.Content.Find.Execute(
    "disegnare la curva di Bézier a partire dai punti di controllo inseriti, utilizzando l’algoritmo di de Casteljau",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "disegnare la curva di Bézier a partire dai punti di controllo inseriti, utilizzando l’algoritmo di de Casteljau",
    2) | Out-Null

# The replace above causes the neighbouring "”. " run to merge with the following
# "Per questa prima richiesta..." run (both share identical formatting). That merge is
# not part of the intended edit, so split them back apart by toggling a formatting
# property on/off across the "Per questa..." span (forces the run boundary to reappear
# without altering the visible formatting).
$full = $d.Content.Text
$idx = $full.IndexOf("Per questa prima richiesta, è stato sufficiente applicare la ")
$len = "Per questa prima richiesta, è stato sufficiente applicare la ".Length
$rsplit = $d.Range($idx, $idx + $len)
$rsplit.Font.Bold = 1
$rsplit.Font.Bold = 0

# --- Change 2: mark the inline image (wp14:anchorId 03A7B561) run as NoProofing ---
$ishp = $d.InlineShapes.Item(2)
$ishp.Range.NoProofing = -1

# --- Change 3: split the trailing run so that ", venga esaminata la " becomes its own run ---
$rng = $d.Content
$rng.Find.Execute("iterazione, ") | Out-Null
$commaStart = $rng.Start + 10
$afterRange = $d.Range($commaStart, $rng.End)
$afterRange.Font.Bold = 1
$afterRange.Text = ", venga esaminata la "
$afterRange.Font.Bold = 0
